# The "reserva_total.prn" sheet (first sheet in the workbook) had its
# previously-generated report contents wiped out as part of verifying the
# system power reserve ("Verificada la potencia del sistema"): all cell
# data, cell formatting, and merged cells on that sheet are removed,
# leaving a completely blank worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Break apart the merged ranges first so Clear() does not leave stray
# mergeCells definitions behind once the cells themselves are empty.
$ws.Cells.UnMerge()

# Remove all cell values/formulas and formatting from the sheet, reverting
# it to a pristine, empty state (dimension collapses back down to A1:A1).
$ws.Cells.Clear()
